# Add a new question row ("Binary Tree Maximum Path Sum" / leetcode 124)
# to the end of the question list on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row is row 25 (row 24 was the previous last data row).
$newRow = 25
$num    = 24   # sequential number shown in column A

$ws.Cells.Item($newRow, 1).Value = $num
$ws.Cells.Item($newRow, 2).Value = "Binary Tree Maximum Path Sum"
$ws.Cells.Item($newRow, 4).Value = "Tree"
$ws.Cells.Item($newRow, 5).Value = "medium"
$ws.Cells.Item($newRow, 6).Value = "leetcode 124"

# Match styling used by the other data rows: column A/D/E/F centered,
# column B left-aligned.
$ws.Cells.Item($newRow, 1).HorizontalAlignment = -4108  # xlCenter
$ws.Cells.Item($newRow, 2).HorizontalAlignment = -4131  # xlLeft
$ws.Cells.Item($newRow, 4).HorizontalAlignment = -4108  # xlCenter
$ws.Cells.Item($newRow, 5).HorizontalAlignment = -4108  # xlCenter
$ws.Cells.Item($newRow, 6).HorizontalAlignment = -4108  # xlCenter

# Update the stored UI selection to match the post-edit state.
$ws.Range("B30").Select()
